# Update the action_code_email_recipient reference values (rows 40-42)
# with the new Malay-language text, and move the active selection
# to C46 (reflecting where the user was last working).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C40").Value = "To: Pengadu, CC: Kumpulan yang berkaitan BCC: Tidak Berkenaan"
$ws.Range("C41").Value = "To: Juruteknik yang berkaitan, CC: Kumpulan Juruteknik yang berkaitan BC: Pengadu"
$ws.Range("C42").Value = "To: Juruteknik yang berkaitan, CC: Kumpulan Juruteknik yang berkaitan BCC: Tidak Berkenaan"

$ws.Range("C46").Select()
